# Update countries & provincias Spain
# - Update "last updated" timestamp in A1
# - Gambia's stats were refreshed, moving it up the ranking (sorted desc by
#   "Casos totales") from just below "San Vicente y las Granadinas" to just
#   below "Laos". This pushes Fiyi, Santa Lucia, Nueva Caledonia, Belice,
#   Islas Virgenes de los Estados Unidos and San Vicente y las Granadinas
#   down by one row each (rows 189-195), while Laos (188) and Namibia (196)
#   stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp: 04:04 -> 04:34
$ws.Range("A1").Value = "Datos actualizados a 8 de Mayo de 2020 a las 04:34"

# New data for the affected block of rows (188-196), column order:
# Pais, Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$data = @(
    @("Laos",                                  19, 0, 9,  10, 0, 0, 0),
    @("Gambia",                                18, 0, 9,  8,  0, 0, 1),
    @("Fiyi",                                  18, 0, 14, 4,  0, 0, 0),
    @("Santa Lucia",                           18, 0, 17, 1,  0, 0, 0),
    @("Nueva Caledonia",                       18, 0, 18, 0,  0, 0, 0),
    @("Belice",                                18, 0, 16, 0,  0, 0, 2),
    @("Islas Virgenes de los Estados Unidos",  17, 0, 0,  17, 0, 0, 0),
    @("San Vicente y las Granadinas",          17, 0, 9,  8,  0, 0, 0),
    @("Namibia",                               16, 0, 9,  7,  0, 0, 0)
)

$startRow = 188
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
    $ws.Cells.Item($row, 7).Value = $values[6]
    $ws.Cells.Item($row, 8).Value = $values[7]
}
